# Auto-generated: update leve profit figures (columns H-N) across all 8 profession sheets
# per the upstream market-data refresh (scheduled runner commit).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 41.25
$ws.Cells.Item(8, 9).Value = 41.25
$ws.Cells.Item(8, 11).Value = 123.75
$ws.Cells.Item(8, 13).Value = 15.25
$ws.Cells.Item(70, 8).Value = 1002.5122
$ws.Cells.Item(70, 9).Value = 1143.75
$ws.Cells.Item(70, 10).Value = 968.2727
$ws.Cells.Item(70, 11).Value = 3431.25
$ws.Cells.Item(70, 12).Value = 2904.8181
$ws.Cells.Item(70, 13).Value = -3161.25
$ws.Cells.Item(70, 14).Value = -3444.8181
$ws.Cells.Item(73, 8).Value = 1002.5122
$ws.Cells.Item(73, 9).Value = 1143.75
$ws.Cells.Item(73, 10).Value = 968.2727
$ws.Cells.Item(73, 11).Value = 3431.25
$ws.Cells.Item(73, 12).Value = 2904.8181
$ws.Cells.Item(73, 13).Value = -2495.25
$ws.Cells.Item(73, 14).Value = -4776.8181
$ws.Cells.Item(97, 8).Value = 1683.6364
$ws.Cells.Item(97, 10).Value = 1527.5
$ws.Cells.Item(97, 12).Value = 4582.5
$ws.Cells.Item(97, 14).Value = -5574.5
$ws.Cells.Item(112, 8).Value = 2977.7297
$ws.Cells.Item(112, 9).Value = 1380
$ws.Cells.Item(112, 10).Value = 3227.375
$ws.Cells.Item(112, 11).Value = 4140
$ws.Cells.Item(112, 12).Value = 9682.125
$ws.Cells.Item(112, 13).Value = -3032
$ws.Cells.Item(112, 14).Value = -11898.125
$ws.Cells.Item(137, 8).Value = 2157.3215
$ws.Cells.Item(137, 9).Value = 1535.25
$ws.Cells.Item(137, 11).Value = 4605.75
$ws.Cells.Item(137, 13).Value = -2055.75
$ws.Cells.Item(141, 8).Value = 2782.0908
$ws.Cells.Item(141, 9).Value = 2733.6191
$ws.Cells.Item(141, 10).Value = 3800
$ws.Cells.Item(141, 11).Value = 8200.8573
$ws.Cells.Item(141, 12).Value = 11400
$ws.Cells.Item(141, 13).Value = -3020.8573
$ws.Cells.Item(141, 14).Value = -21760

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2508.9443
$ws.Cells.Item(2, 9).Value = 1769.1818
$ws.Cells.Item(2, 10).Value = 3671.4285
$ws.Cells.Item(2, 11).Value = 1769.1818
$ws.Cells.Item(2, 12).Value = 3671.4285
$ws.Cells.Item(2, 13).Value = -1656.1818
$ws.Cells.Item(2, 14).Value = -3897.4285
$ws.Cells.Item(116, 8).Value = 2508.9443
$ws.Cells.Item(116, 9).Value = 1769.1818
$ws.Cells.Item(116, 10).Value = 3671.4285
$ws.Cells.Item(116, 11).Value = 1769.1818
$ws.Cells.Item(116, 12).Value = 3671.4285
$ws.Cells.Item(116, 13).Value = 524.8181999999999
$ws.Cells.Item(116, 14).Value = -8259.4285
$ws.Cells.Item(122, 8).Value = 2037.0588
$ws.Cells.Item(122, 9).Value = 1951.5385
$ws.Cells.Item(122, 10).Value = 2315
$ws.Cells.Item(122, 11).Value = 5854.6155
$ws.Cells.Item(122, 12).Value = 6945
$ws.Cells.Item(122, 13).Value = -3404.6155
$ws.Cells.Item(122, 14).Value = -11845

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2508.9443
$ws.Cells.Item(3, 9).Value = 1769.1818
$ws.Cells.Item(3, 10).Value = 3671.4285
$ws.Cells.Item(3, 11).Value = 1769.1818
$ws.Cells.Item(3, 12).Value = 3671.4285
$ws.Cells.Item(3, 13).Value = -1655.1818
$ws.Cells.Item(3, 14).Value = -3899.4285
$ws.Cells.Item(82, 8).Value = 20235
$ws.Cells.Item(82, 9).Value = 7338.6665
$ws.Cells.Item(82, 10).Value = 24786.646
$ws.Cells.Item(82, 11).Value = 7338.6665
$ws.Cells.Item(82, 12).Value = 24786.646
$ws.Cells.Item(82, 13).Value = -6955.6665
$ws.Cells.Item(82, 14).Value = -25552.646
$ws.Cells.Item(85, 8).Value = 20235
$ws.Cells.Item(85, 9).Value = 7338.6665
$ws.Cells.Item(85, 10).Value = 24786.646
$ws.Cells.Item(85, 11).Value = 7338.6665
$ws.Cells.Item(85, 12).Value = 24786.646
$ws.Cells.Item(85, 13).Value = -6012.6665
$ws.Cells.Item(85, 14).Value = -27438.646

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(50, 8).Value = 8959.333000000001
$ws.Cells.Item(50, 10).Value = 8959.333000000001
$ws.Cells.Item(50, 12).Value = 8959.333000000001
$ws.Cells.Item(50, 14).Value = -10209.333
$ws.Cells.Item(51, 8).Value = 7865.25
$ws.Cells.Item(51, 10).Value = 9457
$ws.Cells.Item(51, 12).Value = 9457
$ws.Cells.Item(51, 14).Value = -10929
$ws.Cells.Item(53, 8).Value = 24340
$ws.Cells.Item(53, 10).Value = 24340
$ws.Cells.Item(53, 12).Value = 24340
$ws.Cells.Item(53, 14).Value = -25554
$ws.Cells.Item(60, 8).Value = 18198.223
$ws.Cells.Item(60, 10).Value = 22514
$ws.Cells.Item(60, 12).Value = 22514
$ws.Cells.Item(60, 14).Value = -23536
$ws.Cells.Item(61, 8).Value = 7865.25
$ws.Cells.Item(61, 10).Value = 9457
$ws.Cells.Item(61, 12).Value = 9457
$ws.Cells.Item(61, 14).Value = -10153
$ws.Cells.Item(109, 8).Value = 10857.143
$ws.Cells.Item(109, 10).Value = 10857.143
$ws.Cells.Item(109, 12).Value = 10857.143
$ws.Cells.Item(109, 14).Value = -12937.143
$ws.Cells.Item(122, 8).Value = 1571.7307
$ws.Cells.Item(122, 9).Value = 1056.0526
$ws.Cells.Item(122, 10).Value = 2971.4285
$ws.Cells.Item(122, 11).Value = 3168.1578
$ws.Cells.Item(122, 12).Value = 8914.2855
$ws.Cells.Item(122, 13).Value = -718.1578
$ws.Cells.Item(122, 14).Value = -13814.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(70, 8).Value = 7519.5386
$ws.Cells.Item(70, 9).Value = 8603
$ws.Cells.Item(70, 10).Value = 7038
$ws.Cells.Item(70, 11).Value = 25809
$ws.Cells.Item(70, 12).Value = 21114
$ws.Cells.Item(70, 13).Value = -25494
$ws.Cells.Item(70, 14).Value = -21744
$ws.Cells.Item(73, 8).Value = 7519.5386
$ws.Cells.Item(73, 9).Value = 8603
$ws.Cells.Item(73, 10).Value = 7038
$ws.Cells.Item(73, 11).Value = 25809
$ws.Cells.Item(73, 12).Value = 21114
$ws.Cells.Item(73, 13).Value = -24717
$ws.Cells.Item(73, 14).Value = -23298

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 26880.213
$ws.Cells.Item(70, 9).Value = 30671.75
$ws.Cells.Item(70, 10).Value = 5214.2856
$ws.Cells.Item(70, 11).Value = 30671.75
$ws.Cells.Item(70, 12).Value = 5214.2856
$ws.Cells.Item(70, 13).Value = -30401.75
$ws.Cells.Item(70, 14).Value = -5754.2856
$ws.Cells.Item(73, 8).Value = 26880.213
$ws.Cells.Item(73, 9).Value = 30671.75
$ws.Cells.Item(73, 10).Value = 5214.2856
$ws.Cells.Item(73, 11).Value = 30671.75
$ws.Cells.Item(73, 12).Value = 5214.2856
$ws.Cells.Item(73, 13).Value = -29735.75
$ws.Cells.Item(73, 14).Value = -7086.2856
$ws.Cells.Item(80, 8).Value = 3268.4443
$ws.Cells.Item(80, 9).Value = 2940
$ws.Cells.Item(80, 10).Value = 3362.2856
$ws.Cells.Item(80, 11).Value = 2940
$ws.Cells.Item(80, 12).Value = 3362.2856
$ws.Cells.Item(80, 13).Value = -1942
$ws.Cells.Item(80, 14).Value = -5358.2856
$ws.Cells.Item(83, 8).Value = 3268.4443
$ws.Cells.Item(83, 9).Value = 2940
$ws.Cells.Item(83, 10).Value = 3362.2856
$ws.Cells.Item(83, 11).Value = 14700
$ws.Cells.Item(83, 12).Value = 16811.428
$ws.Cells.Item(83, 13).Value = -9708
$ws.Cells.Item(83, 14).Value = -26795.428
$ws.Cells.Item(122, 8).Value = 1543.5
$ws.Cells.Item(122, 9).Value = 1516.6154
$ws.Cells.Item(122, 10).Value = 1660
$ws.Cells.Item(122, 11).Value = 4549.8462
$ws.Cells.Item(122, 12).Value = 4980
$ws.Cells.Item(122, 13).Value = -2099.8462
$ws.Cells.Item(122, 14).Value = -9880

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3090.45
$ws.Cells.Item(122, 9).Value = 2371.5557
$ws.Cells.Item(122, 10).Value = 3678.6365
$ws.Cells.Item(122, 11).Value = 7114.6671
$ws.Cells.Item(122, 12).Value = 11035.9095
$ws.Cells.Item(122, 13).Value = -4664.6671
$ws.Cells.Item(122, 14).Value = -15935.9095
$ws.Cells.Item(136, 8).Value = 2660.147
$ws.Cells.Item(136, 9).Value = 1865.6818
$ws.Cells.Item(136, 11).Value = 5597.0454
$ws.Cells.Item(136, 13).Value = -3047.0454

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(109, 8).Value = 29088.5
$ws.Cells.Item(109, 10).Value = 29088.5
$ws.Cells.Item(109, 12).Value = 29088.5
$ws.Cells.Item(109, 14).Value = -31862.5
